$wb = $excel.ActiveWorkbook

# Sheets that get the two new accessory rows (MX-DPBX, MX-BBX) inserted just
# above the existing "Wg" / "Accessories" footer rows. $lastRow is the row
# number that currently holds "Wg" (t="s" v="0") on each of these sheets;
# "Accessories" (v="9") is always the row right after it.
$targets = @(
    @{ Name = "Slovakia";    LastRow = 11 },
    @{ Name = "Italy";       LastRow = 11 },
    @{ Name = "Netherlands"; LastRow = 11 },
    @{ Name = "Austria";     LastRow = 10 },
    @{ Name = "Denmark";     LastRow = 11 }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Name)
    $wgRow = $t.LastRow
    $accRow = $wgRow + 1
    $newRow1 = $wgRow
    $newRow2 = $wgRow + 1
    $movedWgRow = $wgRow + 2
    $movedAccRow = $wgRow + 3

    # Remember the existing footer values before they get overwritten.
    $wgValue = $ws.Cells.Item($wgRow, 1).Value()
    $accValue = $ws.Cells.Item($accRow, 1).Value()

    # Push the footer rows down by two, carrying the same cell style.
    $ws.Cells.Item($wgRow, 1).Copy() | Out-Null
    $ws.Cells.Item($movedWgRow, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($accRow, 1).Copy() | Out-Null
    $ws.Cells.Item($movedAccRow, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Cells.Item($movedWgRow, 1).Value = $wgValue
    $ws.Cells.Item($movedAccRow, 1).Value = $accValue

    # Give the two freshly inserted rows the same style as the old footer row.
    $ws.Cells.Item($wgRow, 1).Copy() | Out-Null
    $ws.Cells.Item($newRow1, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($newRow2, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    # Populate the new accessory rows.
    $ws.Cells.Item($newRow1, 1).Value = "MX-DPBX"
    $ws.Cells.Item($newRow2, 1).Value = "MX-BBX"

    $ws.Range($ws.Cells.Item($newRow1, 1), $ws.Cells.Item($newRow2, 1)).Select() | Out-Null
}

# Spain was the active tab before; Denmark becomes the active tab now.
$wb.Worksheets.Item("Denmark").Activate()
